# Update filtered_output.xlsx: refresh the filtered stock data on both
# sheets (Neg_Change / Pos_Change) with the newly computed rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Neg_Change" (sheet 1): 13 existing rows get new values and a
# 14th data row is appended (dimension grows from A1:I13 to A1:I14).
# ---------------------------------------------------------------------
$wsNeg = $wb.Worksheets.Item("Neg_Change")

$negData = @(
    ,@("AXISBANK", 1268, 1273, 1257.8, 1270.1, 4101288, 10175969, -0.596963394837386, "AXISBANK")
    ,@("HINDALCO", 797.15, 798.8, 790, 792, 3005726, 7299555, -0.5882316113790498, "HINDALCO")
    ,@("HINDZINC", 475, 478.5, 473.7, 477.9, 2595527, 6175474, -0.5797040032878448, "HINDZINC")
    ,@("VEDL", 509, 514.8, 507.55, 511.7, 4009545, 7919977, -0.4937428479905939, "VEDL")
    ,@("GAIL", 184.4, 184.88, 183.3, 184.05, 3032584, 6258414, -0.5154388955412665, "GAIL")
    ,@("ALKEM", 5730, 5747.5, 5659.5, 5743, 104658, 243124, -0.5695283065431631, "ALKEM")
    ,@("NYKAA", 269, 270.47, 265.53, 268.99, 3887160, 8683313, -0.5523413701659724, "NYKAA")
    ,@("AUROPHARMA", 1233.9, 1244.3, 1228.7, 1232.7, 744634, 1827771, -0.5925999482429691, "AUROPHARMA")
    ,@("CONCOR", 525.4, 528.4, 516, 517, 1589480, 3630058, -0.5621337179736522, "CONCOR")
    ,@("LUPIN", 2055, 2055, 2019.7, 2025, 684024, 1357360, -0.4960629457181588, "LUPIN")
    ,@("BIOCON", 424, 424.45, 408.45, 409.8, 3249725, 7436143, -0.5629824493692496, "BIOCON")
    ,@("IIFL", 565.9, 565.9, 556.4, 558, 864945, 2048417, -0.57774955001838, "IIFL")
    ,@("LAURUSLABS", 1010, 1010, 982.15, 985, 1288996, 2922214, -0.5588974660993342, "LAURUSLABS")
)

for ($i = 0; $i -lt $negData.Count; $i++) {
    $row = $negData[$i]
    $r = $i + 2
    $wsNeg.Cells.Item($r, 1).Value2 = $row[0]
    $wsNeg.Cells.Item($r, 2).Value2 = $row[1]
    $wsNeg.Cells.Item($r, 3).Value2 = $row[2]
    $wsNeg.Cells.Item($r, 4).Value2 = $row[3]
    $wsNeg.Cells.Item($r, 5).Value2 = $row[4]
    $wsNeg.Cells.Item($r, 6).Value2 = $row[5]
    $wsNeg.Cells.Item($r, 7).Value2 = $row[6]
    $wsNeg.Cells.Item($r, 8).Value2 = $row[7]
    $wsNeg.Cells.Item($r, 9).Value2 = $row[8]
}

# ---------------------------------------------------------------------
# Sheet "Pos_Change" (sheet 2): only the first 10 rows survive (with new
# values); rows 12-16 (old rows 12 through 16) are removed so the
# dimension shrinks from A1:I16 to A1:I11.
# ---------------------------------------------------------------------
$wsPos = $wb.Worksheets.Item("Pos_Change")

$posData = @(
    ,@("INFY", 1490.5, 1543, 1490.5, 1542, 10981126, 7525119, 0.4592627704625003, "INFY")
    ,@("SUNPHARMA", 1762, 1790.3, 1751.1, 1786.5, 2131124, 1429522, 0.4907948251233629, "SUNPHARMA")
    ,@("APOLLOHOSP", 7419, 7482, 7332, 7450, 454714, 285148, 0.5946596153576389, "APOLLOHOSP")
    ,@("CIPLA", 1511.4, 1529.6, 1505.4, 1526, 1171078, 747913, 0.5657944172651097, "CIPLA")
    ,@("NAUKRI", 1324.9, 1364, 1318, 1363.2, 862715, 561417, 0.5366741655489592, "NAUKRI")
    ,@("ADANIGREEN", 1080.2, 1084.5, 1052.2, 1069.9, 2155114, 1414894, 0.5231628659108032, "ADANIGREEN")
    ,@("HAVELLS", 1465, 1466.2, 1438.3, 1444.3, 1021976, 708532, 0.4423851004612354, "HAVELLS")
    ,@("JSWENERGY", 526.05, 528.7, 512.35, 513, 2913764, 2011169, 0.448791225401744, "JSWENERGY")
    ,@("OBEROIRLTY", 1701.9, 1722.7, 1689.1, 1714.8, 334415, 233440, 0.4325522618231665, "OBEROIRLTY")
    ,@("KFINTECH", 1091.9, 1093.3, 1072.2, 1082.9, 461199, 322950, 0.4280817464003716, "KFINTECH")
)

for ($i = 0; $i -lt $posData.Count; $i++) {
    $row = $posData[$i]
    $r = $i + 2
    $wsPos.Cells.Item($r, 1).Value2 = $row[0]
    $wsPos.Cells.Item($r, 2).Value2 = $row[1]
    $wsPos.Cells.Item($r, 3).Value2 = $row[2]
    $wsPos.Cells.Item($r, 4).Value2 = $row[3]
    $wsPos.Cells.Item($r, 5).Value2 = $row[4]
    $wsPos.Cells.Item($r, 6).Value2 = $row[5]
    $wsPos.Cells.Item($r, 7).Value2 = $row[6]
    $wsPos.Cells.Item($r, 8).Value2 = $row[7]
    $wsPos.Cells.Item($r, 9).Value2 = $row[8]
}

# Remove the now-obsolete trailing rows (old rows 12 through 16).
$wsPos.Range("A12:I16").EntireRow.Delete()
